# Update the "想去人数" (want-to-go count) figures in column F for the
# sheets "展览" (sheet1) and "全部类型" (sheet4), leaving "演出" and
# "本地生活" untouched, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 474
$ws1.Range("F4").Value  = 7936
$ws1.Range("F10").Value = 463
$ws1.Range("F13").Value = 450
$ws1.Range("F17").Value = 5826
$ws1.Range("F19").Value = 256
$ws1.Range("F20").Value = 1706
$ws1.Range("F22").Value = 379

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 474
$ws4.Range("F4").Value  = 7936
$ws4.Range("F10").Value = 463
$ws4.Range("F13").Value = 450
$ws4.Range("F18").Value = 5826
$ws4.Range("F21").Value = 256
$ws4.Range("F22").Value = 1706
$ws4.Range("F24").Value = 379
